# Reverse the order of the comma-separated "Recorded By" names/emails
# stored in column G for every data row on the active sheet.
#   "System, dnasr281@gmail.com"        -> "dnasr281@gmail.com, System"
#   "System, system, backup@backdoor.com" -> "backup@backdoor.com, system, System"
# Cells holding a single value (or no value) are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ($text -ne $null -and $text -ne "") {
        $parts = $text -split ", "

        if ($parts.Count -gt 1) {
            $reversed = $parts[($parts.Count - 1)..0]
            $cell.Value = $reversed -join ", "
        }
    }
}
